$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-02 Thursday" "2024-05-03 Friday"

Replace-Text "762×4=" "310×6="
Replace-Text "638×5=" "311×6="
Replace-Text "598×2=" "176×8="
Replace-Text "447×7=" "288×7="
Replace-Text "974×8=" "222×5="

Replace-Text "569×7=" "257×2="
Replace-Text "712×5=" "716×7="
Replace-Text "681×4=" "740×5="
Replace-Text "652×7=" "873×3="
Replace-Text "170×5=" "392×7="

Replace-Text "563×2=" "493×9="
Replace-Text "507×4=" "150×9="
Replace-Text "562×2=" "509×6="
Replace-Text "650×2=" "229×2="
Replace-Text "659×9=" "598×4="

Replace-Text "707×3=" "772×4="
Replace-Text "518×3=" "288×5="
Replace-Text "490×8=" "845×5="
Replace-Text "470×3=" "266×4="
Replace-Text "300×5=" "538×8="

Replace-Text "744×8=" "410×4="
Replace-Text "753×3=" "894×7="
Replace-Text "168×8=" "227×2="
Replace-Text "659×7=" "676×6="
Replace-Text "125×4=" "141×8="
